# Auto-generated Excel COM-interop script
# Applies updated market/profit figures ("chore: update Sheets via scheduled runner")
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("K33").Value = 259.08334
$ws.Range("I33").Value = 259.08334
$ws.Range("H33").Value = 259.08334
$ws.Range("M33").Value = -30.08334000000002
$ws.Range("H100").Value = 2253.2222
$ws.Range("M100").Value = -1404
$ws.Range("K100").Value = 1945
$ws.Range("N100").Value = -3581.8
$ws.Range("L100").Value = 2499.8
$ws.Range("I100").Value = 1945
$ws.Range("J100").Value = 2499.8
$ws.Range("H106").Value = 12823868
$ws.Range("M106").Value = -83334549
$ws.Range("K106").Value = 83335180
$ws.Range("I106").Value = 83335180
$ws.Range("N106").Value = -4890.0908
$ws.Range("L106").Value = 3628.0908
$ws.Range("J106").Value = 3628.0908
$ws.Range("I116").Value = 2750
$ws.Range("J116").Value = 4346.5386
$ws.Range("H116").Value = 4133.6665
$ws.Range("M116").Value = 692
$ws.Range("K116").Value = 2750
$ws.Range("N116").Value = -11230.5386
$ws.Range("L116").Value = 4346.5386

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 11118.333
$ws.Range("N24").Value = -11866.333
$ws.Range("L24").Value = 11118.333
$ws.Range("J24").Value = 11118.333
$ws.Range("H32").Value = 5894.289
$ws.Range("N32").Value = -12686.071
$ws.Range("L32").Value = 12112.071
$ws.Range("J32").Value = 12112.071
$ws.Range("M63").Value = -1279.1
$ws.Range("L63").Value = 15625500
$ws.Range("K63").Value = 1965.1
$ws.Range("N63").Value = -15626872
$ws.Range("J63").Value = 15625500
$ws.Range("I63").Value = 1965.1
$ws.Range("H63").Value = 2605887.5
$ws.Range("M66").Value = -6393.5
$ws.Range("K66").Value = 9825.5
$ws.Range("L66").Value = 78127500
$ws.Range("N66").Value = -78134364
$ws.Range("I66").Value = 1965.1
$ws.Range("J66").Value = 15625500
$ws.Range("H66").Value = 2605887.5
$ws.Range("N88").Value = -501839
$ws.Range("L88").Value = 501027
$ws.Range("I88").Value = 2333.3333
$ws.Range("J88").Value = 501027
$ws.Range("H88").Value = 201810.8
$ws.Range("M88").Value = -1927.3333
$ws.Range("K88").Value = 2333.3333
$ws.Range("M91").Value = -929.3332999999998
$ws.Range("L91").Value = 501027
$ws.Range("K91").Value = 2333.3333
$ws.Range("N91").Value = -503835
$ws.Range("I91").Value = 2333.3333
$ws.Range("J91").Value = 501027
$ws.Range("H91").Value = 201810.8
$ws.Range("H100").Value = 11118.333
$ws.Range("N100").Value = -13282.333
$ws.Range("L100").Value = 11118.333
$ws.Range("J100").Value = 11118.333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H97").Value = 16572.8
$ws.Range("M97").Value = -4241
$ws.Range("K97").Value = 5232
$ws.Range("L97").Value = 24133.334
$ws.Range("N97").Value = -26115.334
$ws.Range("J97").Value = 24133.334
$ws.Range("I97").Value = 5232
$ws.Range("J99").Value = 0
$ws.Range("I99").Value = 2250
$ws.Range("H99").Value = 2250
$ws.Range("K99").Value = 2250
$ws.Range("M99").Value = -752
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()
$ws.Range("I105").Value = 1670.5834
$ws.Range("H105").Value = 1614727.4
$ws.Range("M105").Value = 76.41660000000002
$ws.Range("K105").Value = 1670.5834
$ws.Range("I134").Value = 3889.8386
$ws.Range("H134").Value = 4191.4863
$ws.Range("M134").Value = -9134.515800000001
$ws.Range("K134").Value = 11669.5158

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("L62").Value = 5166.6665
$ws.Range("J62").Value = 5166.6665
$ws.Range("I62").Value = 5483.3335
$ws.Range("H62").Value = 5377.778
$ws.Range("M62").Value = -4859.3335
$ws.Range("K62").Value = 5483.3335
$ws.Range("N62").Value = -6414.6665
$ws.Range("M65").Value = -24296.6675
$ws.Range("K65").Value = 27416.6675
$ws.Range("N65").Value = -32073.3325
$ws.Range("L65").Value = 25833.3325
$ws.Range("I65").Value = 5483.3335
$ws.Range("J65").Value = 5166.6665
$ws.Range("H65").Value = 5377.778
$ws.Range("J68").Value = 64990
$ws.Range("H68").Value = 64990
$ws.Range("L68").Value = 64990
$ws.Range("N68").Value = -66488
$ws.Range("N71").Value = -202458
$ws.Range("L71").Value = 194970
$ws.Range("J71").Value = 64990
$ws.Range("H71").Value = 64990
$ws.Range("L74").Value = 47557.25
$ws.Range("N74").Value = -49305.25
$ws.Range("J74").Value = 47557.25
$ws.Range("H74").Value = 47557.25
$ws.Range("N77").Value = -151407.75
$ws.Range("L77").Value = 142671.75
$ws.Range("J77").Value = 47557.25
$ws.Range("H77").Value = 47557.25
$ws.Range("M87").Value = -3814
$ws.Range("K87").Value = 5000
$ws.Range("L87").Value = 10000
$ws.Range("N87").Value = -12372
$ws.Range("I87").Value = 5000
$ws.Range("J87").Value = 10000
$ws.Range("H87").Value = 7500
$ws.Range("L90").Value = 30000
$ws.Range("I90").Value = 5000
$ws.Range("J90").Value = 10000
$ws.Range("H90").Value = 7500
$ws.Range("M90").Value = -9072
$ws.Range("K90").Value = 15000
$ws.Range("N90").Value = -41856
$ws.Range("H97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
$ws.Range("J97").Value = 0
$ws.Range("L111").Value = 43000
$ws.Range("N111").Value = -51180
$ws.Range("J111").Value = 43000
$ws.Range("H111").Value = 43000
$ws.Range("L132").Value = 18805.9995
$ws.Range("N132").Value = -23865.9995
$ws.Range("I132").Value = 1458.9474
$ws.Range("J132").Value = 6268.6665
$ws.Range("H132").Value = 2613.28
$ws.Range("M132").Value = -1846.8422
$ws.Range("K132").Value = 4376.8422

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("L131").Value = 2234.6739
$ws.Range("N131").Value = -12314.6739
$ws.Range("J131").Value = 744.8913
$ws.Range("H131").Value = 717.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("J15").Value = 18750
$ws.Range("H15").Value = 18750
$ws.Range("N15").Value = -19326
$ws.Range("L15").Value = 18750
$ws.Range("L81").Value = 18750
$ws.Range("J81").Value = 18750
$ws.Range("H81").Value = 18750
$ws.Range("N81").Value = -20746
$ws.Range("H84").Value = 18750
$ws.Range("N84").Value = -66234
$ws.Range("L84").Value = 56250
$ws.Range("J84").Value = 18750
$ws.Range("H97").Value = 1383.25
$ws.Range("M97").Value = -917.4375
$ws.Range("K97").Value = 1413.4375
$ws.Range("L97").Value = 1262.5
$ws.Range("N97").Value = -2254.5
$ws.Range("J97").Value = 1262.5
$ws.Range("I97").Value = 1413.4375
$ws.Range("I99").Value = 0
$ws.Range("H99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("L7").Value = 6350
$ws.Range("I7").Value = 2829.5881
$ws.Range("J7").Value = 6350
$ws.Range("H7").Value = 3200.158
$ws.Range("M7").Value = -2717.5881
$ws.Range("K7").Value = 2829.5881
$ws.Range("N7").Value = -6574
$ws.Range("I68").Value = 933.3333
$ws.Range("J68").Value = 2713.5715
$ws.Range("H68").Value = 2179.5
$ws.Range("M68").Value = -184.3333
$ws.Range("K68").Value = 933.3333
$ws.Range("L68").Value = 2713.5715
$ws.Range("N68").Value = -4211.5715
$ws.Range("K71").Value = 4666.6665
$ws.Range("N71").Value = -21055.8575
$ws.Range("L71").Value = 13567.8575
$ws.Range("J71").Value = 2713.5715
$ws.Range("I71").Value = 933.3333
$ws.Range("H71").Value = 2179.5
$ws.Range("M71").Value = -922.6665000000003
$ws.Range("N126").Value = -23990
$ws.Range("L126").Value = 19050
$ws.Range("J126").Value = 6350
$ws.Range("I126").Value = 2829.5881
$ws.Range("H126").Value = 3200.158
$ws.Range("M126").Value = -6018.764299999999
$ws.Range("K126").Value = 8488.764299999999
$ws.Range("I132").Value = 1337.7778
$ws.Range("H132").Value = 2077.4167
$ws.Range("M132").Value = -1483.3334
$ws.Range("K132").Value = 4013.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 47619630
$ws.Range("M107").Value = -249998880
$ws.Range("K107").Value = 250000800
$ws.Range("N107").Value = -6869.3334
$ws.Range("I107").Value = 83333600
$ws.Range("L107").Value = 3029.3334
$ws.Range("J107").Value = 1009.7778
$ws.Range("I122").Value = 1702.9
$ws.Range("J122").Value = 1850
$ws.Range("H122").Value = 1727.4166
$ws.Range("L122").Value = 5550
$ws.Range("M122").Value = -2658.700000000001
$ws.Range("K122").Value = 5108.700000000001
$ws.Range("N122").Value = -10450
